$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the existing two-row block (rows 11:12) onto the
# two new rows (13:14) so the new entries pick up the same cell styles
# (date column, time columns, duration formula column, notes columns).
$ws.Range("A11:F12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 13: Source MUX architecture ---
$ws.Range("A13").Value = "30.3.2020"
$ws.Range("B13").Value = 0.49652777777777773
$ws.Range("C13").Value = 0.50347222222222221
$ws.Range("E13").Value = "Source MUX"
$ws.Range("F13").Value = "Create necessary files"

# --- Row 14: Source MUX testbench ---
$ws.Range("A14").Value = "30.3.2020"
$ws.Range("B14").Value = 0.50347222222222221
$ws.Range("C14").Value = 0.52083333333333337
$ws.Range("E14").Value = "Source MUX"
$ws.Range("F14").Value = "RTL and TB"

# Extend the "To - From" duration formula down through the two new rows as
# one shared formula block, matching how the column was filled for the
# earlier rows.
$ws.Range("D13:D14").Formula = "=C13-B13"

# Match the cursor position left behind after entering the new data.
$ws.Range("F15").Select() | Out-Null
